# Update the F1 score shown in the metrics table on slide 16 from 0.947 to 0.95
# ("add makedown to 90%" commit: the f1 row of the Precision/Recall/Accuracy/f1
# table is rounded to match the other values already shown in that table).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# Find the specific table cell that currently reads "0.947" and update it,
# scanning every table on the slide so the exact shape index doesn't matter.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cellRange = $tbl.Cell($r, $c).Shape.TextFrame.TextRange
                if ($cellRange.Text -eq "0.947") {
                    $cellRange.Text = "0.95"
                }
            }
        }
    }
}
